$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Fitness) values for rows 43-135 per the run_6 log edit.
$ws.Range("C43:C46").Value = 7320
$ws.Range("C47:C74").Value = 7295
$ws.Range("C75:C135").Value = 7293
